$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "10 2015"
